# Updates the crypto price table (columns D "Price" and E "Volume(1h)")
# for rows 2-51, matching the latest scrape snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '70.615.10'
$ws.Range("E2").Value = '  -0.99%  '

$ws.Range("D3").Value = '3.789.76'
$ws.Range("E3").Value = '  -2.04%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  -0.07%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '709.77'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +2.55%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '169.90'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -1.97%  '

$ws.Range("D7").Value = '3.787.18'
$ws.Range("E7").Value = '  -2.00%  '

$ws.Range("E8").Value = '  +0.02%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.521'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -1.27%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.160'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -2.19%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.36'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -0.31%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.456'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -1.69%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000253'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -2.42%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '36.16'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -1.43%  '

$ws.Range("D15").Value = '4.426.95'
$ws.Range("E15").Value = '  -2.16%  '

$ws.Range("D16").Value = '3.849.72'
$ws.Range("E16").Value = '  -0.87%  '

$ws.Range("D17").Value = '70.608.24'
$ws.Range("E17").Value = '  -1.03%  '

$ws.Range("E18").Value = '  +0.00%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.13'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -2.14%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.33'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -2.89%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '492.36'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -0.89%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.59'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -4.46%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.725'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +0.05%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '84.82'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -0.18%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000145'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -1.86%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.06'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -2.70%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.47'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -1.33%  '

$ws.Range("D28").Value = '3.936.99'
$ws.Range("E28").Value = '  -2.31%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +0.03%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.05'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -4.97%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.10'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -1.07%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.31'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -4.63%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.22'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -3.99%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '29.08'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -2.58%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.175'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -1.46%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.998'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -0.21%  '

$ws.Range("D37").Value = '3.755.79'
$ws.Range("E37").Value = '  -1.72%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '9.02'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -3.41%  '

$ws.Range("E39").Value = '  -3.10%  '

$ws.Range("E40").Value = '  +1.42%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.31'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -3.55%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.91'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -2.63%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.28'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -4.83%  '

$ws.Range("E45").Value = '  +0.00%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '164.35'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -0.04%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.000311'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +1.35%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '48.94'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +0.51%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '420.57'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +1.35%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.66'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -0.57%  '

$ws.Range("E51").Value = '  -1.49%  '
